$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 261, shifting existing rows 261-402 down to 262-403
$ws.Rows.Item(261).Insert()

# Populate the new row. Values are set in C, B, A order so that the
# newly-created shared strings are appended to sharedStrings.xml in the
# same order they appear in the canonical workbook (long note, tax id,
# company name -> indices 1091, 1092, 1093).
$ws.Range("C261").Value = "95437539 尋光有限公司`r`n臺北市中山區中原里新生北路２段３１之１號１１樓之６`r`n（未向國際貿易署登記出進口廠商資料者，出口金額限制美金兩萬以下，且通關必驗，若金額超過美金兩萬需檢附輸出許可證才可出口）"
$ws.Range("B261").Value = "95437539"
$ws.Range("A261").Value = "尋光有限公司"

# Match the row height used by similarly-structured rows (e.g. row 260).
$ws.Rows.Item(261).RowHeight = 78

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the
# (now one-row-larger) data range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=" + $ws.Name + "!`$A`$1:`$C`$402"
    }
}

# Update the view so the active cell follows the last data row, which
# shifted from row 402 to row 403.
$excel.Goto($ws.Range("A389"), $true)
$ws.Range("A403").Select()
